# Rearranged pin location settings
# - Flip sign of several pin-position values on Sheet1
# - Update Sheet1's view (selection)
# - Build out Sheet2 with "Conventions" / "Examples" documentation content

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: fix the sign of a handful of mirrored pin-position values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("G4").Value  = -1.9690000000000001

$ws1.Range("F7").Value  = -4.8
$ws1.Range("G7").Value  = -4.8

$ws1.Range("F10").Value = 1.9690000000000001
$ws1.Range("G10").Value = 1.9690000000000001

$ws1.Range("F13").Value = 4.8
$ws1.Range("G13").Value = 4.8

# ---------------------------------------------------------------------------
# Sheet2: populate the "Conventions" / "Examples" notes
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Column widths (character units)
$ws2.Columns.Item(2).ColumnWidth = 47.833333333333336
$ws2.Columns.Item(3).ColumnWidth = 18.666666666666668

# Row 1 - section title
$ws2.Range("A1").Value = "Conventions"

# Row 2 - "b or p" convention explanation
$ws2.Range("A2").Value = "b or p"
$ws2.Range("B2").Value = "base or patform pin position"

# Row 3 - "bAB" convention explanation (first char regular, rest smaller rich-text run)
$ws2.Range("A3").Value = "bAB"
$ws2.Range("A3").Characters(2, 2).Font.Size = 10

# Row 5 - Examples title
$ws2.Range("A5").Value = "Examples"

# Row 6 - base pin example (reuses existing shared string "b11")
$ws2.Range("A6").Value = "b11"
$ws2.Range("B6").Value = "Base pin, X1"

# Row 7 - platform pin example (reuses existing shared string "p36")
$ws2.Range("A7").Value = "p36"
$ws2.Range("B7").Value = "Platform pin, Z3"

# Row 3, column B - wrapped explanatory text (added last so it becomes the
# final new shared string)
$ws2.Range("B3").Value = "A is the dof x,y ,or z .  B is the actuator number X1, X2, Y1, Z1, Z2, Z3"
$ws2.Range("B3").WrapText = $true
$ws2.Range("B3").Font.Size = 10
$ws2.Range("B3").Font.Size = 11

$ws2.Rows.Item(3).RowHeight = 30

$ws2.PageSetup.Orientation = 1

# View state
$ws2.Range("B6").Select() | Out-Null

$ws1.Select()
$ws1.Range("F38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Edit applied"
